$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number but must remain text
# (matching the original inline-string cell type): force text format first,
# then clear the format back to default after writing so no extra cell
# style lingers on the sheet.
$textForceCells = @("D5", "D8", "D16", "D18", "D22", "D25", "D27", "D36", "D37", "D41", "D42", "D47", "D50")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.176.22"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "1.682.41"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "215.14"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "22.68"
$ws.Range("E8").Value = "  +5.00%  "
$ws.Range("E9").Value = "  +2.18%  "
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("D12").Value = "1.921.10"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").Value = "1.683.41"
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("E14").Value = "  +2.11%  "
$ws.Range("E15").Value = "  +4.81%  "
$ws.Range("D16").Value = "66.81"
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("D17").Value = "27.166.83"
$ws.Range("E17").Value = "  +0.54%  "
$ws.Range("D18").Value = "236.21"
$ws.Range("E18").Value = "  +0.21%  "
$ws.Range("E19").Value = "  -3.31%  "
$ws.Range("D20").Value = "0.0₃0740"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").Value = "4.54"
$ws.Range("E22").Value = "  +1.83%  "
$ws.Range("E24").Value = "  -1.49%  "
$ws.Range("D25").Value = "146.92"
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("E26").Value = "  +2.31%  "
$ws.Range("D27").Value = "16.33"
$ws.Range("E27").Value = "  -1.73%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  +1.02%  "
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("E32").Value = "  +0.37%  "
$ws.Range("D33").Value = "1.547.62"
$ws.Range("E33").Value = "  +1.88%  "
$ws.Range("E34").Value = "  +2.15%  "
$ws.Range("E35").Value = "  -2.57%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "0.950"
$ws.Range("E36").Value = "  +3.52%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "0.605"
$ws.Range("E37").Value = "  +2.57%  "
$ws.Range("E38").Value = "  -0.24%  "
$ws.Range("E39").Value = "  -1.31%  "
$ws.Range("E40").Value = "  +3.48%  "
$ws.Range("D41").Value = "5.78"
$ws.Range("E41").Value = "  +1.01%  "
$ws.Range("D42").Value = "69.20"
$ws.Range("E42").Value = "  +1.79%  "
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("D45").Value = "1.828.39"
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("E46").Value = "  +1.43%  "
$ws.Range("D47").Value = "90.09"
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("E48").Value = "  +4.15%  "
$ws.Range("E49").Value = "  +6.65%  "
$ws.Range("D50").Value = "8.20"
$ws.Range("E50").Value = "  +3.43%  "
$ws.Range("E51").Value = "  +0.16%  "

foreach ($addr in $textForceCells) {
    $ws.Range($addr).ClearFormats()
}
